$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Font colours for the two "prediction" blocks.
#    Order matters: the workbook's styles.xml lists the RED font before the
#    GREEN font, so apply RED first (creates font index 2) then GREEN
#    (creates font index 3) to reproduce the same font/cellXf ordering.
#    Excel's Font.Color is a BGR-packed long, not a hex RGB string:
#      FFFF0000 (red)   -> R=0xFF G=0x00 B=0x00 -> 255
#      FF00B050 (green) -> R=0x00 G=0xB0 B=0x50 -> 5287936
# ---------------------------------------------------------------------------
$ws.Range("A36:A65").Font.Color = 255
$ws.Range("A3:A33").Font.Color = 5287936

# ---------------------------------------------------------------------------
# 2) New cells added in column B (same text already used elsewhere in the
#    sheet) plus the corrected A45 value. Use Value2 on the read side since
#    Value's getter is unreliable for round-tripping through this host; the
#    Value setter itself works fine.
# ---------------------------------------------------------------------------
$ws.Range("B18").Value = $ws.Range("B17").Value2
$ws.Range("B33").Value = $ws.Range("B30").Value2
$ws.Range("A45").Value = $ws.Range("B47").Value2
$ws.Range("B50").Value = $ws.Range("B47").Value2
$ws.Range("B65").Value = $ws.Range("B62").Value2

# ---------------------------------------------------------------------------
# 3) View state: scroll down and move the selection to A28.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("A28").Select()
